$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.586.98'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '2.033.94'
$ws.Range("E3").Value = '  +2.88%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.601'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.15'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0752'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("D13").Value = '2.333.09'
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.760'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("E17").Value = '  +1.97%  '
$ws.Range("D18").Value = '2.064.71'
$ws.Range("E18").Value = '  +4.54%  '
$ws.Range("D19").Value = '36.739.58'
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +17.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").Value = '0.0₃0797'
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '220.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.14%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.126'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("E31").Value = '  +5.60%  '
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0607'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.31%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.75'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.65%  '
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +35.63%  '
$ws.Range("D43").Value = '1.476.08'
$ws.Range("E43").Value = '  +3.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0942'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.10%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("E51").Value = '  +4.92%  '
